$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4702473333333333
$ws.Range("N2").Value = 1.410742
$ws.Range("O2").Value = 0.009034922268422819
$ws.Range("P2").Value = 0.009034922268422819
$ws.Range("Q2").Value = 0.6532774706606667
$ws.Range("R2").Value = 5.879497235946
$ws.Range("S2").Value = 0.0026294067397609
$ws.Range("T2").Value = 0.0026294067397609
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("N3").Value = 0.9584440000000001
$ws.Range("O3").Value = 0.006138235792679485
$ws.Range("P3").Value = 0.006138235792679485
$ws.Range("Q3").Value = 0.4438301773746667
$ws.Range("R3").Value = 3.994471596372001
$ws.Range("S3").Value = 0.00178639263117097
$ws.Range("T3").Value = 0.00178639263117097
$ws.Range("G4").Value = 1.389221
$ws.Range("H4").Value = 4.167663
$ws.Range("I4").Value = 0.2910270461264192
$ws.Range("J4").Value = 0.2910270461264192
$ws.Range("M4").Value = 1.047307
$ws.Range("N4").Value = 3.141921
$ws.Range("O4").Value = 0.02012204358311108
$ws.Range("P4").Value = 0.02012204358311108
$ws.Range("Q4").Value = 1.454940877847
$ws.Range("R4").Value = 13.094467900623
$ws.Range("S4").Value = 0.005856058906019886
$ws.Range("T4").Value = 0.005856058906019886
$ws.Range("G5").Value = 1.389221
$ws.Range("H5").Value = 4.167663
$ws.Range("I5").Value = 0.2910270461264192
$ws.Range("J5").Value = 0.2910270461264192
$ws.Range("M5").Value = 50.21070966666667
$ws.Range("N5").Value = 150.632129
$ws.Range("O5").Value = 0.9647047983557866
$ws.Range("P5").Value = 0.9647047983557866
$ws.Range("Q5").Value = 69.75377229383635
$ws.Range("R5").Value = 627.7839506445271
$ws.Range("S5").Value = 0.2807551878494674
$ws.Range("T5").Value = 0.2807551878494674
$ws.Range("I6").Value = 0.461328155686921
$ws.Range("J6").Value = 0.4613281556869209
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4702473333333333
$ws.Range("N6").Value = 1.410742
$ws.Range("O6").Value = 0.009034922268422819
$ws.Range("P6").Value = 0.009034922268422819
$ws.Range("Q6").Value = 1.035557673085778
$ws.Range("R6").Value = 9.320019057772001
$ws.Range("S6").Value = 0.004168064026866192
$ws.Range("T6").Value = 0.004168064026866192
$ws.Range("I7").Value = 0.461328155686921
$ws.Range("J7").Value = 0.4613281556869209
$ws.Range("N7").Value = 0.9584440000000001
$ws.Range("O7").Value = 0.006138235792679485
$ws.Range("P7").Value = 0.006138235792679485
$ws.Range("Q7").Value = 0.7035475221004446
$ws.Range("R7").Value = 6.331927698904002
$ws.Range("S7").Value = 0.002831740997408272
$ws.Range("T7").Value = 0.002831740997408272
$ws.Range("I8").Value = 0.461328155686921
$ws.Range("J8").Value = 0.4613281556869209
$ws.Range("M8").Value = 1.047307
$ws.Range("N8").Value = 3.141921
$ws.Range("O8").Value = 0.02012204358311108
$ws.Range("P8").Value = 0.02012204358311108
$ws.Range("Q8").Value = 2.306332695687334
$ws.Range("R8").Value = 20.756994261186
$ws.Range("S8").Value = 0.009282865254848479
$ws.Range("T8").Value = 0.009282865254848479
$ws.Range("I9").Value = 0.461328155686921
$ws.Range("J9").Value = 0.4613281556869209
$ws.Range("M9").Value = 50.21070966666667
$ws.Range("N9").Value = 150.632129
$ws.Range("O9").Value = 0.9647047983557866
$ws.Range("P9").Value = 0.9647047983557866
$ws.Range("Q9").Value = 110.5717820829016
$ws.Range("R9").Value = 995.1460387461143
$ws.Range("S9").Value = 0.445045485407798
$ws.Range("T9").Value = 0.445045485407798
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1506176666666667
$ws.Range("H10").Value = 0.451853
$ws.Range("I10").Value = 0.03155280162368235
$ws.Range("J10").Value = 0.03155280162368235
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4702473333333333
$ws.Range("N10").Value = 1.410742
$ws.Range("O10").Value = 0.009034922268422819
$ws.Range("P10").Value = 0.009034922268422819
$ws.Range("Q10").Value = 0.07082755610288889
$ws.Range("R10").Value = 0.637448004926
$ws.Range("S10").Value = 0.0002850771100209354
$ws.Range("T10").Value = 0.0002850771100209354
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1506176666666667
$ws.Range("H11").Value = 0.451853
$ws.Range("I11").Value = 0.03155280162368235
$ws.Range("J11").Value = 0.03155280162368235
$ws.Range("N11").Value = 0.9584440000000001
$ws.Range("O11").Value = 0.006138235792679485
$ws.Range("P11").Value = 0.006138235792679485
$ws.Range("Q11").Value = 0.04811953297022223
$ws.Range("R11").Value = 0.433075796732
$ws.Range("S11").Value = 0.0001936785362858024
$ws.Range("T11").Value = 0.0001936785362858024
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1506176666666667
$ws.Range("H12").Value = 0.451853
$ws.Range("I12").Value = 0.03155280162368235
$ws.Range("J12").Value = 0.03155280162368235
$ws.Range("M12").Value = 1.047307
$ws.Range("N12").Value = 3.141921
$ws.Range("O12").Value = 0.02012204358311108
$ws.Range("P12").Value = 0.02012204358311108
$ws.Range("Q12").Value = 0.1577429366236667
$ws.Range("R12").Value = 1.419686429613
$ws.Range("S12").Value = 0.0006349068494409945
$ws.Range("T12").Value = 0.0006349068494409945
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1506176666666667
$ws.Range("H13").Value = 0.451853
$ws.Range("I13").Value = 0.03155280162368235
$ws.Range("J13").Value = 0.03155280162368235
$ws.Range("M13").Value = 50.21070966666667
$ws.Range("N13").Value = 150.632129
$ws.Range("O13").Value = 0.9647047983557866
$ws.Range("P13").Value = 0.9647047983557866
$ws.Range("Q13").Value = 7.562619931670779
$ws.Range("R13").Value = 68.06357938503702
$ws.Range("S13").Value = 0.03043913912793462
$ws.Range("T13").Value = 0.03043913912793462
$ws.Range("G14").Value = 1.031517666666667
$ws.Range("H14").Value = 3.094553
$ws.Range("I14").Value = 0.2160919965629775
$ws.Range("J14").Value = 0.2160919965629775
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4702473333333333
$ws.Range("N14").Value = 1.410742
$ws.Range("O14").Value = 0.009034922268422819
$ws.Range("P14").Value = 0.009034922268422819
$ws.Range("Q14").Value = 0.4850684320362221
$ws.Range("R14").Value = 4.365615888325999
$ws.Range("S14").Value = 0.001952374391774793
$ws.Range("T14").Value = 0.001952374391774793
$ws.Range("G15").Value = 1.031517666666667
$ws.Range("H15").Value = 3.094553
$ws.Range("I15").Value = 0.2160919965629775
$ws.Range("J15").Value = 0.2160919965629775
$ws.Range("N15").Value = 0.9584440000000001
$ws.Range("O15").Value = 0.006138235792679485
$ws.Range("P15").Value = 0.006138235792679485
$ws.Range("Q15").Value = 0.3295506395035555
$ws.Range("R15").Value = 2.965955755532
$ws.Range("S15").Value = 0.001326423627814441
$ws.Range("T15").Value = 0.001326423627814441
$ws.Range("G16").Value = 1.031517666666667
$ws.Range("H16").Value = 3.094553
$ws.Range("I16").Value = 0.2160919965629775
$ws.Range("J16").Value = 0.2160919965629775
$ws.Range("M16").Value = 1.047307
$ws.Range("N16").Value = 3.141921
$ws.Range("O16").Value = 0.02012204358311108
$ws.Range("P16").Value = 0.02012204358311108
$ws.Range("Q16").Value = 1.080315672923666
$ws.Range("R16").Value = 9.722841056312999
$ws.Range("S16").Value = 0.004348212572801724
$ws.Range("T16").Value = 0.004348212572801724
$ws.Range("G17").Value = 1.031517666666667
$ws.Range("H17").Value = 3.094553
$ws.Range("I17").Value = 0.2160919965629775
$ws.Range("J17").Value = 0.2160919965629775
$ws.Range("M17").Value = 50.21070966666667
$ws.Range("N17").Value = 150.632129
$ws.Range("O17").Value = 0.9647047983557866
$ws.Range("P17").Value = 0.9647047983557866
$ws.Range("Q17").Value = 51.79323407703745
$ws.Range("R17").Value = 466.1391066933371
$ws.Range("S17").Value = 0.2084649859705866
$ws.Range("T17").Value = 0.2084649859705866
